$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4398.8  # was 4499
$ws.Range("I113").Value = 4373.5  # was 4498.5
$ws.Range("K113").Value = 4373.5  # was 4498.5
$ws.Range("M113").Value = -1119.5  # was -1244.5
$ws.Range("H116").Value = 11627.5  # was 14337
$ws.Range("J116").Value = 5501.6665  # was 6503
$ws.Range("L116").Value = 5501.6665  # was 6503
$ws.Range("N116").Value = -12385.6665  # was -13387
$ws.Range("H135").Value = 1500.9286  # was 1434.2
$ws.Range("J135").Value = 7265.5  # was 5010.3335
$ws.Range("L135").Value = 65389.5  # was 45093.0015
$ws.Range("N135").Value = -70459.5  # was -50163.0015
$ws.Range("H137").Value = 1394.5834  # was 0
$ws.Range("I137").Value = 1273.7  # was 0
$ws.Range("J137").Value = 1999  # was 0
$ws.Range("K137").Value = 3821.1  # was 0
$ws.Range("L137").Value = 5997  # was 0
$ws.Range("M137").Value = -1271.1  # was None
$ws.Range("N137").Value = -11097  # was None
$ws.Range("H138").Value = 3404.5  # was 3673.1052
$ws.Range("I138").Value = 2733.2222  # was 3287.25
$ws.Range("K138").Value = 8199.6666  # was 9861.75
$ws.Range("M138").Value = -3059.6666  # was -4721.75
$ws.Range("H141").Value = 2287  # was 2731.8333
$ws.Range("I141").Value = 1793.909  # was 2278.4
$ws.Range("K141").Value = 5381.727000000001  # was 6835.200000000001
$ws.Range("M141").Value = -201.7270000000008  # was -1655.200000000001

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2266915.5  # was 2335609.2
$ws.Range("I32").Value = 2123488.8  # was 2189847
$ws.Range("K32").Value = 2123488.8  # was 2189847
$ws.Range("M32").Value = -2123201.8  # was -2189560
$ws.Range("H62").Value = 100000  # was 0
$ws.Range("J62").Value = 100000  # was 0
$ws.Range("L62").Value = 100000  # was 0
$ws.Range("N62").Value = -101248  # was None
$ws.Range("H63").Value = 4061.5  # was 3665.6667
$ws.Range("I63").Value = 4141.7144  # was 3999.1428
$ws.Range("J63").Value = 3500  # was 2498.5
$ws.Range("K63").Value = 4141.7144  # was 3999.1428
$ws.Range("L63").Value = 3500  # was 2498.5
$ws.Range("M63").Value = -3455.7144  # was -3313.1428
$ws.Range("N63").Value = -4872  # was -3870.5
$ws.Range("H65").Value = 100000  # was 0
$ws.Range("J65").Value = 100000  # was 0
$ws.Range("L65").Value = 300000  # was 0
$ws.Range("N65").Value = -306240  # was None
$ws.Range("H66").Value = 4061.5  # was 3665.6667
$ws.Range("I66").Value = 4141.7144  # was 3999.1428
$ws.Range("J66").Value = 3500  # was 2498.5
$ws.Range("K66").Value = 20708.572  # was 19995.714
$ws.Range("L66").Value = 17500  # was 12492.5
$ws.Range("M66").Value = -17276.572  # was -16563.714
$ws.Range("N66").Value = -24364  # was -19356.5
$ws.Range("H102").Value = 3563.5  # was 4418
$ws.Range("I102").Value = 3084.6667  # was 4127
$ws.Range("K102").Value = 3084.6667  # was 4127
$ws.Range("M102").Value = -1462.6667  # was -2505

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 700  # was 449.5
$ws.Range("I22").Value = 0  # was 199
$ws.Range("K22").Value = 0  # was 199
$ws.Range("M22").ClearContents()  # was -26
$ws.Range("H94").Value = 398.6  # was 400.1
$ws.Range("I94").Value = 425.125  # was 389
$ws.Range("J94").Value = 292.5  # was 500
$ws.Range("K94").Value = 425.125  # was 389
$ws.Range("L94").Value = 292.5  # was 500
$ws.Range("M94").Value = 25.875  # was 62
$ws.Range("N94").Value = -1194.5  # was -1402

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2478.875  # was 2375.5557
$ws.Range("I31").Value = 2261.5715  # was 2172.5
$ws.Range("K31").Value = 2261.5715  # was 2172.5
$ws.Range("M31").Value = -1966.5715  # was -1877.5
$ws.Range("H34").Value = 2478.875  # was 2375.5557
$ws.Range("I34").Value = 2261.5715  # was 2172.5
$ws.Range("K34").Value = 2261.5715  # was 2172.5
$ws.Range("M34").Value = -2059.5715  # was -1970.5
$ws.Range("H58").Value = 1741.6666  # was 1786.3636
$ws.Range("I58").Value = 1555.3  # was 1589.2222
$ws.Range("K58").Value = 1555.3  # was 1589.2222
$ws.Range("M58").Value = -1352.3  # was -1386.2222
$ws.Range("H86").Value = 9050.362999999999  # was 9815.444
$ws.Range("I86").Value = 9356.799999999999  # was 10294.125
$ws.Range("K86").Value = 9356.799999999999  # was 10294.125
$ws.Range("M86").Value = -8233.799999999999  # was -9171.125
$ws.Range("H89").Value = 9050.362999999999  # was 9815.444
$ws.Range("I89").Value = 9356.799999999999  # was 10294.125
$ws.Range("K89").Value = 46784  # was 51470.625
$ws.Range("M89").Value = -41168  # was -45854.625
$ws.Range("H122").Value = 3209.0908  # was 3050
$ws.Range("I122").Value = 1500  # was 1466.6666
$ws.Range("K122").Value = 4500  # was 4399.9998
$ws.Range("M122").Value = -2050  # was -1949.9998
$ws.Range("H134").Value = 3396.6843  # was 3545.3125
$ws.Range("I134").Value = 2910.5  # was 2951.4
$ws.Range("J134").Value = 4230.143  # was 4535.1665
$ws.Range("K134").Value = 8731.5  # was 8854.200000000001
$ws.Range("L134").Value = 12690.429  # was 13605.4995
$ws.Range("M134").Value = -6196.5  # was -6319.200000000001
$ws.Range("N134").Value = -17760.429  # was -18675.4995
$ws.Range("H136").Value = 1741.6666  # was 1786.3636
$ws.Range("I136").Value = 1555.3  # was 1589.2222
$ws.Range("K136").Value = 4665.9  # was 4767.6666
$ws.Range("M136").Value = -2115.9  # was -2217.6666

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 239.72728  # was 183.3
$ws.Range("I12").Value = 352.33334  # was 217.2
$ws.Range("J12").Value = 104.6  # was 149.4
$ws.Range("K12").Value = 1057.00002  # was 651.5999999999999
$ws.Range("L12").Value = 313.8  # was 448.2
$ws.Range("M12").Value = -884.0000199999999  # was -478.5999999999999
$ws.Range("N12").Value = -659.8  # was -794.2
$ws.Range("H22").Value = 2337.5  # was 2000
$ws.Range("I22").Value = 2000  # was 0
$ws.Range("J22").Value = 2368.182  # was 2000
$ws.Range("K22").Value = 6000  # was 0
$ws.Range("L22").Value = 7104.545999999999  # was 6000
$ws.Range("M22").Value = -5831  # was None
$ws.Range("N22").Value = -7442.545999999999  # was -6338
$ws.Range("H27").Value = 2337.5  # was 2000
$ws.Range("I27").Value = 2000  # was 0
$ws.Range("J27").Value = 2368.182  # was 2000
$ws.Range("K27").Value = 6000  # was 0
$ws.Range("L27").Value = 7104.545999999999  # was 6000
$ws.Range("M27").Value = -5898  # was None
$ws.Range("N27").Value = -7308.545999999999  # was -6204

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 999  # was 999.5
$ws.Range("I22").Value = 999  # was 999.5
$ws.Range("K22").Value = 999  # was 999.5
$ws.Range("M22").Value = -470  # was -470.5
$ws.Range("H70").Value = 3335.6667  # was 3001.75
$ws.Range("I70").Value = 3335.6667  # was 3001.75
$ws.Range("K70").Value = 3335.6667  # was 3001.75
$ws.Range("M70").Value = -3065.6667  # was -2731.75
$ws.Range("H73").Value = 3335.6667  # was 3001.75
$ws.Range("I73").Value = 3335.6667  # was 3001.75
$ws.Range("K73").Value = 3335.6667  # was 3001.75
$ws.Range("M73").Value = -2399.6667  # was -2065.75
$ws.Range("H80").Value = 2568.2222  # was 2775.889
$ws.Range("I80").Value = 2247  # was 2436.8
$ws.Range("J80").Value = 2969.75  # was 3199.75
$ws.Range("K80").Value = 2247  # was 2436.8
$ws.Range("L80").Value = 2969.75  # was 3199.75
$ws.Range("M80").Value = -1249  # was -1438.8
$ws.Range("N80").Value = -4965.75  # was -5195.75
$ws.Range("H83").Value = 2568.2222  # was 2775.889
$ws.Range("I83").Value = 2247  # was 2436.8
$ws.Range("J83").Value = 2969.75  # was 3199.75
$ws.Range("K83").Value = 11235  # was 12184
$ws.Range("L83").Value = 14848.75  # was 15998.75
$ws.Range("M83").Value = -6243  # was -7192
$ws.Range("N83").Value = -24832.75  # was -25982.75

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7849.1113  # was 7849.278
$ws.Range("I7").Value = 9058.200000000001  # was 9058.799999999999
$ws.Range("K7").Value = 9058.200000000001  # was 9058.799999999999
$ws.Range("M7").Value = -8946.200000000001  # was -8946.799999999999
$ws.Range("H17").Value = 0  # was 5000
$ws.Range("J17").Value = 0  # was 5000
$ws.Range("L17").Value = 0  # was 5000
$ws.Range("N17").ClearContents()  # was -5340
$ws.Range("H82").Value = 1049.5  # was 1027.8572
$ws.Range("I82").Value = 649.5  # was 749.5
$ws.Range("J82").Value = 1249.5  # was 1139.2
$ws.Range("K82").Value = 649.5  # was 749.5
$ws.Range("L82").Value = 1249.5  # was 1139.2
$ws.Range("M82").Value = -288.5  # was -388.5
$ws.Range("N82").Value = -1971.5  # was -1861.2
$ws.Range("H85").Value = 1049.5  # was 1027.8572
$ws.Range("I85").Value = 649.5  # was 749.5
$ws.Range("J85").Value = 1249.5  # was 1139.2
$ws.Range("K85").Value = 649.5  # was 749.5
$ws.Range("L85").Value = 1249.5  # was 1139.2
$ws.Range("M85").Value = 598.5  # was 498.5
$ws.Range("N85").Value = -3745.5  # was -3635.2
$ws.Range("H122").Value = 6216.76  # was 6108.423
$ws.Range("I122").Value = 4985.7856  # was 4880.067
$ws.Range("K122").Value = 14957.3568  # was 14640.201
$ws.Range("M122").Value = -12507.3568  # was -12190.201
$ws.Range("H126").Value = 7849.1113  # was 7849.278
$ws.Range("I126").Value = 9058.200000000001  # was 9058.799999999999
$ws.Range("K126").Value = 27174.6  # was 27176.4
$ws.Range("M126").Value = -24704.6  # was -24706.4
$ws.Range("H136").Value = 2999.125  # was 3160.2222
$ws.Range("I136").Value = 2999.125  # was 3160.2222
$ws.Range("K136").Value = 8997.375  # was 9480.6666
$ws.Range("M136").Value = -6447.375  # was -6930.6666

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2500001.2  # was 2000001.2
$ws.Range("H84").Value = 2500001.2  # was 2000001.2
$ws.Range("H95").Value = 0  # was 14500
$ws.Range("J95").Value = 0  # was 14500
$ws.Range("L95").Value = 0  # was 14500
$ws.Range("N95").ClearContents()  # was -19992
$ws.Range("H107").Value = 4047.9167  # was 4122.0835
$ws.Range("I107").Value = 3997.4443  # was 4096.3335
$ws.Range("K107").Value = 11992.3329  # was 12289.0005
$ws.Range("M107").Value = -10072.3329  # was -10369.0005
$ws.Range("H136").Value = 1965.625  # was 1895.0588
$ws.Range("I136").Value = 2030  # was 1951
$ws.Range("K136").Value = 6090  # was 5853
$ws.Range("M136").Value = -3540  # was -3303
